$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New user record: "TEST" / "TEST" / "6BF02F00" / counter 28 / roomID 105 / access TRUE ---
# Write the new row 4 values first (leading apostrophe forces text storage for
# values that would otherwise be auto-typed as number/boolean by Excel).
$ws.Range("A4").Value = "TEST"
$ws.Range("B4").Value = "TEST"
$ws.Range("C4").Value = "6BF02F00"
$ws.Range("D4").Value = 28
$ws.Range("E4").Value = "'105"
$ws.Range("F4").Value = "'TRUE"

# Copy the formatting from the old row 3 onto the new row 4 so the cell
# styles line up with the rest of the table.
$ws.Range("A3:F3").Copy()
$ws.Range("A4:F4").PasteSpecial(-4122)  # xlPasteFormats

# Clear out the old row 3 record (user has been migrated to row 4) while
# keeping its formatting in place.
$ws.Range("A3:F3").ClearContents()
